$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so values like "0.9998" or
# "1.0000" are stored as literal strings (matching the source data feed)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '24.988.92'
$ws.Range('E2').Value = '  -3.71%  '

$ws.Range('D3').Value = '1.645.08'
$ws.Range('E3').Value = '  -5.53%  '

$ws.Range('D4').Value = '0.9998'
$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').Value = '235.67'
$ws.Range('E5').Value = '  -5.68%  '

$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').Value = '0.4819'
$ws.Range('E7').Value = '  -6.10%  '

$ws.Range('D8').Value = '0.2592'
$ws.Range('E8').Value = '  -5.72%  '

$ws.Range('D9').Value = '0.05995'
$ws.Range('E9').Value = '  -3.05%  '

$ws.Range('D10').Value = '0.07188'
$ws.Range('E10').Value = '  -0.49%  '

$ws.Range('D11').Value = '1.645.22'
$ws.Range('E11').Value = '  -5.55%  '

$ws.Range('D12').Value = '14.75'
$ws.Range('E12').Value = '  -2.20%  '

$ws.Range('D13').Value = '0.6193'
$ws.Range('E13').Value = '  -4.59%  '

$ws.Range('D14').Value = '4.496'
$ws.Range('E14').Value = '  -2.85%  '

$ws.Range('D15').Value = '72.75'
$ws.Range('E15').Value = '  -6.20%  '

$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('D17').Value = '0.9998'
$ws.Range('E17').Value = '  -0.14%  '

$ws.Range('D18').Value = '24.985.37'
$ws.Range('E18').Value = '  -3.82%  '

$ws.Range('D19').Value = '11.38'
$ws.Range('E19').Value = '  -3.72%  '

$ws.Range('D20').Value = '0.000006608'
$ws.Range('E20').Value = '  -2.72%  '

$ws.Range('D21').Value = '4.493'
$ws.Range('E21').Value = '  +5.34%  '

$ws.Range('D22').Value = '1.857.95'
$ws.Range('E22').Value = '  -5.57%  '

$ws.Range('D23').Value = '8.606'
$ws.Range('E23').Value = '  -0.76%  '

$ws.Range('D24').Value = '5.276'
$ws.Range('E24').Value = '  -1.95%  '

$ws.Range('D25').Value = '131.99'
$ws.Range('E25').Value = '  -2.86%  '

$ws.Range('D26').Value = '14.85'

$ws.Range('D27').Value = '1.396'
$ws.Range('E27').Value = '  -7.38%  '

$ws.Range('D28').Value = '102.89'
$ws.Range('E28').Value = '  -2.87%  '

$ws.Range('D29').Value = '1.664'
$ws.Range('E29').Value = '  -6.34%  '

$ws.Range('D30').Value = '3.737'
$ws.Range('E30').Value = '  -5.29%  '

$ws.Range('D31').Value = '0.07833'
$ws.Range('E31').Value = '  -4.52%  '

$ws.Range('D32').Value = '3.566'
$ws.Range('E32').Value = '  -2.20%  '

$ws.Range('D33').Value = '0.04479'
$ws.Range('E33').Value = '  -4.67%  '

$ws.Range('D34').Value = '0.9997'
$ws.Range('E34').Value = '  -0.06%  '

$ws.Range('D35').Value = '2.591'
$ws.Range('E35').Value = '  -2.55%  '

$ws.Range('D36').Value = '0.9299'
$ws.Range('E36').Value = '  -6.68%  '

$ws.Range('D37').Value = '0.5817'
$ws.Range('E37').Value = '  -6.84%  '

$ws.Range('D38').Value = '2.573'
$ws.Range('E38').Value = '  -5.72%  '

$ws.Range('D39').Value = '0.01566'
$ws.Range('E39').Value = '  -2.96%  '

$ws.Range('D40').Value = '0.8594'
$ws.Range('E40').Value = '  +13.71%  '

$ws.Range('D41').Value = '0.9999'
$ws.Range('E41').Value = '  -0.14%  '

$ws.Range('B42').Value = 'PaxosStandard'
$ws.Range('C42').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D42').Value = '1.001'
$ws.Range('E42').Value = '  -0.16%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.814'
$ws.Range('E43').Value = '  -5.20%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '98.14'
$ws.Range('E44').Value = '  -1.83%  '

$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.3705'
$ws.Range('E45').Value = '  -3.71%  '

$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '4.778'
$ws.Range('E46').Value = '  -4.78%  '

$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1149'
$ws.Range('E47').Value = '  +1.81%  '

$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '6.089'
$ws.Range('E48').Value = '  -3.31%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.05191'
$ws.Range('E49').Value = '  -0.69%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '29.72'
$ws.Range('E50').Value = '  -3.39%  '

$ws.Range('B51').Value = 'TrueUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D51').Value = '1.0000'
$ws.Range('E51').Value = '  -0.35%  '
